$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "39.604.11"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.04%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.171.72"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.66"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.33%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.62%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "63.19"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.392"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.67%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0850"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.22%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.44%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.90"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.85%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.492.88"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.58%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.79"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.20%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.812"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.15%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.35%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.166.72"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.42%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "39.608.99"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0917"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +7.31%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.99%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.49%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "229.45"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -3.53%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.24%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.60"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.00%  "
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "170.93"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.08%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.16%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.46"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.86"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.99%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +3.70%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.40%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.39%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.70"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.88%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.77%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.84"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +7.32%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.48%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.25%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.90"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +16.93%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0229"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.40%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "102.60"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.77%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.73"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.91%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.513.96"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.31%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +2.05%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.73%  "
$ws.Range("B47").Value = "HuobiToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.80"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.08%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0920"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.02%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000197"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +32.72%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "49.34"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +5.92%  "
